# Leia-me.docx edit
# 1) Merge the split "SYSMA" / bookmark / "TRICULA, tem por " runs in the
#    last paragraph into one run with the full finished sentence (this also
#    removes the stray "_GoBack" bookmark that was sitting mid-word).
# 2) Append four new paragraphs of body text.
# 3) Append three new blank paragraphs, re-homing the "_GoBack" bookmark
#    into the middle one of that trailing trio.

$d = $word.ActiveDocument

# --- helper: append a raw <w:p>...</w:p> fragment at the very end of the
# document body via Range.InsertXML (Flat-OPC wrapped WordprocessingML). ---
function Append-ParagraphXml {
    param($Doc, [string]$InnerXml)

    $r = $Doc.Content
    $r.Collapse(0)  # wdCollapseEnd

    $flatOpc = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' + `
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
        '<pkg:xmlData>' + `
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
        '<w:body>' + $InnerXml + '</w:body></w:document>' + `
        '</pkg:xmlData></pkg:part></pkg:package>'

    $r.InsertXML($flatOpc)
}

# --- 1. fix the mangled sentence -----------------------------------------
# The text "SYSMA" + hidden "_GoBack" bookmark + "TRICULA, tem por " is
# really one sentence split across two runs. Find/Replace across the run
# boundary merges it back into a single run and (because the match engulfs
# the bookmark's location) drops the stray bookmark too.
$searchText = "SYSMATRICULA, tem por "
$replaceText = "SYSMATRICULA, tem por finalidade verificar uma serie de matrículas da empresa XPTO."
$d.Content.Find.Execute($searchText, $false, $false, $false, $false, $false, `
    $true, 1, $false, $replaceText, 2) | Out-Null

# --- 2. three new body paragraphs -----------------------------------------
$lang = '<w:rPr><w:lang w:val="pt-BR"/></w:rPr>'

$para1 = "<w:p><w:pPr>$lang</w:pPr><w:r>$lang<w:t>As matrículas armazenadas em lote (arquivo em formato texto) e separadas com quebra de linha, sendo a divisão de matrículas indicada pela mudança de linha, serão verificadas da seguinte maneira:</w:t></w:r></w:p>"
Append-ParagraphXml $d $para1

$para2 = "<w:p><w:pPr>$lang</w:pPr><w:r>$lang<w:t>* As matriculas que não houverem dígito verificador terão este gerado automaticamente.</w:t></w:r></w:p>"
Append-ParagraphXml $d $para2

$para3 = "<w:p><w:pPr>$lang</w:pPr><w:r>$lang<w:t>* As matrículas que possuírem o dígito verificador serão autenticadas, acusando se está correto ou não.</w:t></w:r></w:p>"
Append-ParagraphXml $d $para3

# --- 3. four trailing blank paragraphs -------------------------------------
$blank = "<w:p><w:pPr>$lang</w:pPr></w:p>"
Append-ParagraphXml $d $blank   # blank #1
Append-ParagraphXml $d $blank   # blank #2 (bookmark goes here)
Append-ParagraphXml $d $blank   # blank #3
Append-ParagraphXml $d $blank   # blank #4

# Re-home the "_GoBack" bookmark onto the second of the four trailing
# blank paragraphs. Adding a bookmark with the same name as an existing
# (even hidden) one moves it, removing the old one automatically.
$paraCount = $d.Paragraphs.Count
$targetPara = $d.Paragraphs($paraCount - 2)
$d.Bookmarks.Add("_GoBack", $targetPara.Range) | Out-Null
